$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.522.66'
$ws.Range('E2').Value = '  +1.39%  '

$ws.Range('D3').Value = '2.514.17'
$ws.Range('E3').Value = '  +0.96%  '

$ws.Range('E4').Value = '  +0.04%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '317.74'
$ws.Range('E5').Value = '  +5.51%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '94.36'
$ws.Range('E6').Value = '  +0.06%  '

$ws.Range('E7').Value = '  -0.25%  '

$ws.Range('E8').Value = '  -0.09%  '

$ws.Range('E9').Value = '  +0.28%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.67'
$ws.Range('E10').Value = '  -0.65%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0812'
$ws.Range('E11').Value = '  +1.30%  '

$ws.Range('E12').Value = '  +2.35%  '

$ws.Range('E13').Value = '  -2.38%  '

$ws.Range('D14').Value = '2.901.21'
$ws.Range('E14').Value = '  +1.20%  '

$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.20'
$ws.Range('E15').Value = '  +2.53%  '

$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '2.491.56'
$ws.Range('E16').Value = '  -0.86%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.845'
$ws.Range('E17').Value = '  -0.22%  '

$ws.Range('D18').Value = '42.609.30'
$ws.Range('E18').Value = '  +1.54%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.88'
$ws.Range('E19').Value = '  +1.36%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.66'
$ws.Range('E20').Value = '  +4.92%  '

$ws.Range('E21').Value = '  -0.27%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '69.25'
$ws.Range('E22').Value = '  -1.92%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '250.62'
$ws.Range('E23').Value = '  +1.44%  '

$ws.Range('E24').Value = '  +2.58%  '

$ws.Range('E25').Value = '  +1.77%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.68'
$ws.Range('E26').Value = '  +0.79%  '

$ws.Range('E27').Value = '  +0.29%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.45'
$ws.Range('E28').Value = '  +6.53%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '41.05'
$ws.Range('E29').Value = '  +11.21%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '10.21'
$ws.Range('E30').Value = '  +1.80%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.92'
$ws.Range('E31').Value = '  +1.39%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '157.06'
$ws.Range('E32').Value = '  +2.49%  '

$ws.Range('E33').Value = '  +3.58%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '19.09'
$ws.Range('E34').Value = '  +4.92%  '

$ws.Range('E35').Value = '  -0.04%  '

$ws.Range('E36').Value = '  +0.53%  '

$ws.Range('E37').Value = '  +0.52%  '

$ws.Range('E38').Value = '  -2.06%  '

$ws.Range('E39').Value = '  -0.27%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '23.55'
$ws.Range('E40').Value = '  -0.93%  '

$ws.Range('E41').Value = '  +17.09%  '

$ws.Range('E42').Value = '  +0.48%  '

$ws.Range('E43').Value = '  +2.61%  '

$ws.Range('E44').Value = '  -0.61%  '

$ws.Range('E45').Value = '  -0.83%  '

$ws.Range('D46').Value = '2.013.00'
$ws.Range('E46').Value = '  -1.15%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '85.00'
$ws.Range('E47').Value = '  +2.64%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.90'
$ws.Range('E48').Value = '  +0.44%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '74.56'
$ws.Range('E49').Value = '  +4.62%  '

$ws.Range('D50').Value = '2.756.39'
$ws.Range('E50').Value = '  +0.98%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '101.84'
$ws.Range('E51').Value = '  +2.04%  '
